# Update the cryptocurrency price/volume figures in the "Price" (D) and
# "Volume(1h)" (E) columns. A leading apostrophe is used for values that
# look like plain numbers so Excel keeps storing them as text (matching
# the original inline-string cell type) instead of silently converting
# them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.699.34'
$ws.Range('E2').Value = '  +3.40%  '
$ws.Range('D3').Value = '2.444.41'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''575.74'
$ws.Range('E5').Value = '  +2.51%  '
$ws.Range('D6').Value = '''145.64'
$ws.Range('E6').Value = '  +3.25%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +0.65%  '
$ws.Range('D9').Value = '2.444.12'
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('E10').Value = '  +2.05%  '
$ws.Range('D12').Value = '''5.24'
$ws.Range('E12').Value = '  +1.65%  '
$ws.Range('E13').Value = '  +3.00%  '
$ws.Range('D14').Value = '''28.40'
$ws.Range('E14').Value = '  +9.12%  '
$ws.Range('E15').Value = '  +5.64%  '
$ws.Range('D16').Value = '2.888.60'
$ws.Range('E16').Value = '  +2.11%  '
$ws.Range('D17').Value = '62.581.07'
$ws.Range('D18').Value = '2.445.07'
$ws.Range('E18').Value = '  +1.57%  '
$ws.Range('D19').Value = '''7.76'
$ws.Range('E19').Value = '  -3.74%  '
$ws.Range('D20').Value = '''10.91'
$ws.Range('E20').Value = '  +2.80%  '
$ws.Range('D21').Value = '0.0₆0857'
$ws.Range('E21').Value = '  +212.93%  '
$ws.Range('D22').Value = '''326.66'
$ws.Range('E22').Value = '  +1.15%  '
$ws.Range('E24').Value = '  +11.18%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').Value = '''65.53'
$ws.Range('E26').Value = '  +0.99%  '
$ws.Range('D27').Value = '''645.05'
$ws.Range('E27').Value = '  +15.31%  '
$ws.Range('D28').Value = '''1.16'
$ws.Range('E28').Value = '  +16.00%  '
$ws.Range('E29').Value = '  +5.96%  '
$ws.Range('D30').Value = '0.0₃0977'
$ws.Range('E30').Value = '  +4.86%  '
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('E33').Value = '  +6.93%  '
$ws.Range('E34').Value = '  +3.70%  '
$ws.Range('E35').Value = '  +6.01%  '
$ws.Range('E36').Value = '  +2.35%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('E38').Value = '  +3.23%  '
$ws.Range('E39').Value = '  +6.89%  '
$ws.Range('E41').Value = '  +0.90%  '
$ws.Range('E42').Value = '  +1.84%  '
$ws.Range('E43').Value = '  +8.99%  '
$ws.Range('E44').Value = '  +5.73%  '
$ws.Range('D45').Value = '''42.56'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').Value = '''15.01'
$ws.Range('E47').Value = '  +28.02%  '
$ws.Range('D48').Value = '''144.14'
$ws.Range('E48').Value = '  +1.85%  '
$ws.Range('E49').Value = '  +2.11%  '
$ws.Range('D50').Value = '''20.57'
$ws.Range('E50').Value = '  +7.02%  '
$ws.Range('D51').Value = '''0.605'
$ws.Range('E51').Value = '  +3.14%  '
